$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(3,1,197.1981153333333,591.594346,0.5013718116429122,0.5302871392468994,3,1,4.006170333333333,12.018511,0.05496065106112269,0.06019443438116229,790.0092394376451,7110.083154938806,0.02755572119158903,0.03192033440657176),
    @(3,1,197.1981153333333,591.594346,0.5013718116429122,0.5302871392468994,3,1,32.65736066666667,97.972082,0.4480263330901556,0.490690906813236,6439.969975227597,57959.72977704837,0.2246277742851421,0.2602070772284578),
    @(3,1,197.1981153333333,591.594346,0.5013718116429122,0.5302871392468994,3,1,8.814931666666666,26.444795,0.1209320481029573,0.1324481441462082,1738.287911458785,15644.59120312907,0.06063192004306753,0.0702355474578537),
    @(3,1,197.1981153333333,591.594346,0.5013718116429122,0.5302871392468994,3,1,8.399816666666666,25.19945,0.1152370853912109,0.1262108625158624,1656.428015812189,14907.8521423097,0.0577766262710404,0.06692799722542041),
    @(3,1,197.1981153333333,591.594346,0.5013718116429122,0.5302871392468994,2,1,19.0133305,38.026661,0.2608438823545532,0.1904556521435312,3749.392940809785,22496.35764485871,0.130779769852073,0.1009961829285958),
    @(3,1,131.273506,393.820518,0.3337599622221713,0.3530087082119477,3,1,4.006170333333333,12.018511,0.05496065106112269,0.06019443438116229,525.9040252898553,4733.136227608698,0.01834366482186625,0.02124915952244296),
    @(3,1,131.273506,393.820518,0.3337599622221713,0.3530087082119477,3,1,32.65736066666667,97.972082,0.4480263330901556,0.490690906813236,4287.046231419831,38583.41608277847,0.1495332520067083,0.1732181631454897),
    @(3,1,131.273506,393.820518,0.3337599622221713,0.3530087082119477,3,1,8.814931666666666,26.444795,0.1209320481029573,0.1324481441462082,1157.166985033757,10414.50286530381,0.04036227580629284,0.04675534827012279),
    @(3,1,131.273506,393.820518,0.3337599622221713,0.3530087082119477,3,1,8.399816666666666,25.19945,0.1152370853912109,0.1262108625158624,1102.673383590567,9924.0604523151,0.03846152526676369,0.04455353353904033),
    @(3,1,131.273506,393.820518,0.3337599622221713,0.3530087082119477,2,1,19.0133305,38.026661,0.2608438823545532,0.1904556521435312,2495.946555471733,14975.6793328304,0.08705924432054019,0.06723250373485203),
    @(3,1,0.278468,0.835404,0.0007079986814710624,0.0007488306814808822,3,1,4.006170333333333,12.018511,0.05496065106112269,0.06019443438116229,1.115590240382667,10.040312163444,0.00003891206848406601,0.00004507543931900201),
    @(3,1,0.278468,0.835404,0.0007079986814710624,0.0007488306814808822,3,1,32.65736066666667,97.972082,0.4480263330901556,0.490690906813236,9.094029910125334,81.846269191128,0.0003172020530921451,0.0003674444061454275),
    @(3,1,0.278468,0.835404,0.0007079986814710624,0.0007488306814808822,3,1,8.814931666666666,26.444795,0.1209320481029573,0.1324481441462082,2.454676391353333,22.09208752218,0.00008561973060448888,0.00009918123404188318),
    @(3,1,0.278468,0.835404,0.0007079986814710624,0.0007488306814808822,3,1,8.399816666666666,25.19945,0.1152370853912109,0.1262108625158624,2.339080147533333,21.0517213278,0.00008158770451354557,0.00009451056618804317),
    @(3,1,0.278468,0.835404,0.0007079986814710624,0.0007488306814808822,2,1,19.0133305,38.026661,0.2608438823545532,0.1904556521435312,5.294604117674001,31.767624706044,0.0001846771247768166,0.0001426190357865263),
    @(3,1,0.2270173333333333,0.681052,0.0005771865085793579,0.0006104742535155658,3,1,4.006170333333333,12.018511,0.05496065106112269,0.06019443438116229,0.9094701059524444,8.185230953572,0.00003172254629521779,0.00003674715239463176),
    @(3,1,0.2270173333333333,0.681052,0.0005771865085793579,0.0006104742535155658,3,1,32.65736066666667,97.972082,0.4480263330901556,0.490690906813236,7.413786932251556,66.724082390264,0.0002585947549479193,0.0002995541650436863),
    @(3,1,0.2270173333333333,0.681052,0.0005771865085793579,0.0006104742535155658,3,1,8.814931666666666,26.444795,0.1209320481029573,0.1324481441462082,2.001142280482222,18.01028052434,0.0000698003466198969,0.0000808561819271785),
    @(3,1,0.2270173333333333,0.681052,0.0005771865085793579,0.0006104742535155658,3,1,8.399816666666666,25.19945,0.1152370853912109,0.1262108625158624,1.906903980155555,17.1621358214,0.00006651329097581437,0.00007704848207992681),
    @(3,1,0.2270173333333333,0.681052,0.0005771865085793579,0.0006104742535155658,2,1,19.0133305,38.026661,0.2608438823545532,0.1904556521435312,4.316355587895334,25.898133527372,0.0001505555697405094,0.0001162682720701425),
    @(2,1,64.34001,128.68002,0.1635830409448661,0.1153448476061565,3,1,4.006170333333333,12.018511,0.05496065106112269,0.06019443438116229,257.7570393083701,1546.54223585022,0.008990630432888131,0.00694311786043395),
    @(2,1,64.34001,128.68002,0.1635830409448661,0.1153448476061565,3,1,32.65736066666667,97.972082,0.4480263330901556,0.490690906813236,2101.17491186694,12607.04947120164,0.07328950999026514,0.05659866786809943),
    @(2,1,64.34001,128.68002,0.1635830409448661,0.1153448476061565,3,1,8.814931666666666,26.444795,0.1209320481029573,0.1324481441462082,567.15279158265,3402.9167494959,0.01978243217637259,0.01527721100226263),
    @(2,1,64.34001,128.68002,0.1635830409448661,0.1153448476061565,3,1,8.399816666666666,25.19945,0.1152370853912109,0.1262108625158624,540.4442883315,3242.665729989,0.01885083285791749,0.01455777270313372),
    @(2,1,64.34001,128.68002,0.1635830409448661,0.1153448476061565,2,1,19.0133305,38.026661,0.2608438823545532,0.1904556521435312,1223.317874503305,4893.271498013221,0.04266963548742272,0.02196807817222676),
)

for ($r = 0; $r -lt 25; $r++) {
    for ($c = 0; $c -lt 16; $c++) {
        $ws.Cells.Item($r+2, $c+5).Value = $data[$r][$c]
    }
}
